# Updated cryptos list (price + 1h volume change) data refresh.
# Mirrors the GitHub Actions scheduled update of cryptos.xlsx.
# All assigned values are prefixed with a leading apostrophe so Excel
# keeps them as literal text (quote-prefix) instead of auto-converting
# numeric-looking strings (e.g. '1.00', '21.30', '0.0210') into numbers
# and silently dropping the formatting-significant trailing/leading zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'35.806.42"
$ws.Range('E2').Value = "'  -4.36%  "

# Row 3
$ws.Range('D3').Value = "'1.948.96"
$ws.Range('E3').Value = "'  -4.45%  "

# Row 4
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = "'  +0.34%  "

# Row 5
$ws.Range('D5').Value = "'241.71"
$ws.Range('E5').Value = "'  -4.37%  "

# Row 6
$ws.Range('D6').Value = "'0.618"
$ws.Range('E6').Value = "'  -4.50%  "

# Row 7
$ws.Range('D7').Value = "'60.17"
$ws.Range('E7').Value = "'  -7.62%  "

# Row 8
$ws.Range('E8').Value = "'  +0.21%  "

# Row 9
$ws.Range('E9').Value = "'  -2.28%  "

# Row 10
$ws.Range('D10').Value = "'55.94"
$ws.Range('E10').Value = "'  -5.32%  "

# Row 11
$ws.Range('E11').Value = "'  +3.13%  "

# Row 12
$ws.Range('D12').Value = "'0.102"
$ws.Range('E12').Value = "'  -1.60%  "

# Row 13
$ws.Range('D13').Value = "'0.854"
$ws.Range('E13').Value = "'  -5.90%  "

# Row 14
$ws.Range('D14').Value = "'13.77"
$ws.Range('E14').Value = "'  -8.56%  "

# Row 15
$ws.Range('D15').Value = "'2.242.34"
$ws.Range('E15').Value = "'  -4.05%  "

# Row 16
$ws.Range('D16').Value = "'21.30"
$ws.Range('E16').Value = "'  +2.89%  "

# Row 17
$ws.Range('D17').Value = "'5.34"
$ws.Range('E17').Value = "'  -4.38%  "

# Row 18
$ws.Range('D18').Value = "'1.972.82"
$ws.Range('E18').Value = "'  -3.45%  "

# Row 19
$ws.Range('D19').Value = "'35.760.67"
$ws.Range('E19').Value = "'  -4.18%  "

# Row 20
$ws.Range('D20').Value = "'70.34"

# Row 21
$ws.Range('D21').Value = "'0.0₃0840"
$ws.Range('E21').Value = "'  -3.83%  "

# Row 22
$ws.Range('D22').Value = "'236.61"
$ws.Range('E22').Value = "'  +0.10%  "

# Row 23
$ws.Range('D23').Value = "'5.13"
$ws.Range('E23').Value = "'  -4.16%  "

# Row 24
$ws.Range('E24').Value = "'  -0.28%  "

# Row 25
$ws.Range('D25').Value = "'2.48"
$ws.Range('E25').Value = "'  -9.74%  "

# Row 26
$ws.Range('E26').Value = "'  -3.28%  "

# Row 27
$ws.Range('D27').Value = "'9.62"
$ws.Range('E27').Value = "'  +0.85%  "

# Row 28
$ws.Range('D28').Value = "'158.04"
$ws.Range('E28').Value = "'  -4.66%  "

# Row 29
$ws.Range('D29').Value = "'0.133"
$ws.Range('E29').Value = "'  +19.59%  "

# Row 30
$ws.Range('D30').Value = "'19.49"
$ws.Range('E30').Value = "'  -1.86%  "

# Row 31
$ws.Range('E31').Value = "'  -2.85%  "

# Row 32
$ws.Range('D32').Value = "'4.81"
$ws.Range('E32').Value = "'  -7.68%  "

# Row 33
$ws.Range('E33').Value = "'  -7.78%  "

# Row 34
$ws.Range('E34').Value = "'  -1.06%  "

# Row 35
$ws.Range('D35').Value = "'4.32"
$ws.Range('E35').Value = "'  -8.12%  "

# Row 36
$ws.Range('B36').Value = "'THORChain"
$ws.Range('C36').Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range('D36').Value = "'6.22"
$ws.Range('E36').Value = "'  +4.06%  "

# Row 37
$ws.Range('B37').Value = "'BinanceUSD"
$ws.Range('C37').Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range('D37').Value = "'1.00"
$ws.Range('E37').Value = "'  +0.35%  "

# Row 38
$ws.Range('D38').Value = "'1.84"
$ws.Range('E38').Value = "'  +1.70%  "

# Row 39
$ws.Range('E39').Value = "'  -7.49%  "

# Row 40
$ws.Range('D40').Value = "'3.08"
$ws.Range('E40').Value = "'  +13.05%  "

# Row 41
$ws.Range('E41').Value = "'  -6.70%  "

# Row 42
$ws.Range('E42').Value = "'  -2.57%  "

# Row 43
$ws.Range('D43').Value = "'2.80"
$ws.Range('E43').Value = "'  -4.43%  "

# Row 44
$ws.Range('D44').Value = "'0.0210"
$ws.Range('E44').Value = "'  -4.26%  "

# Row 45
$ws.Range('D45').Value = "'1.08"
$ws.Range('E45').Value = "'  -5.24%  "

# Row 46
$ws.Range('D46').Value = "'91.60"
$ws.Range('E46').Value = "'  -3.93%  "

# Row 47
$ws.Range('D47').Value = "'15.91"
$ws.Range('E47').Value = "'  -6.13%  "

# Row 48
$ws.Range('D48').Value = "'7.47"
$ws.Range('E48').Value = "'  -7.37%  "

# Row 49
$ws.Range('D49').Value = "'1.329.23"
$ws.Range('E49').Value = "'  -6.47%  "

# Row 50
$ws.Range('E50').Value = "'  -7.52%  "

# Row 51
$ws.Range('D51').Value = "'2.139.08"
$ws.Range('E51').Value = "'  -3.71%  "
